$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.653.62"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.021.21"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'235.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.82%  "
$ws.Range("D6").Value = "'0.601"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.88%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'54.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.23%  "
$ws.Range("D9").Value = "'0.370"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "'58.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "2.319.19"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "'14.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "'20.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.91%  "
$ws.Range("D16").Value = "'0.763"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "'5.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "2.019.98"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "37.017.06"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "'67.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.30%  "
$ws.Range("D21").Value = "0.0₃0798"
$ws.Range("E21").Value = "  -4.70%  "
$ws.Range("D22").Value = "'5.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.78%  "
$ws.Range("D23").Value = "'220.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.50%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  -7.45%  "
$ws.Range("D27").Value = "'163.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'8.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").Value = "'1.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "'18.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").Value = "'4.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.68%  "
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("D35").Value = "'2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("D36").Value = "'4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.96%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'3.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").Value = "'5.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.22%  "
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("D42").Value = "1.457.94"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'0.0928"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +40.53%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0204"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.81%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'90.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "'15.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "'6.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
